# Rewriting the load structure of df's --> now based on header files
# Update the "dtype" column (C) header label and its numpy dtype values,
# and move the active selection from B14 to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

$ws.Range("C1").Value = "dType"
$ws.Range("C2").Value = "str" + $nbsp
$ws.Range("C3").Value = "str"
$ws.Range("C4").Value = "numpy.float_" + $nbsp
$ws.Range("C5").Value = "str" + $nbsp
$ws.Range("C6").Value = "int" + $nbsp
$ws.Range("C7").Value = "int" + $nbsp
$ws.Range("C8").Value = "str" + $nbsp
$ws.Range("C9").Value = "str" + $nbsp
$ws.Range("C10").Value = "numpy.float_" + $nbsp
$ws.Range("C11").Value = "str"
$ws.Range("C12").Value = "numpy.float_" + $nbsp
$ws.Range("C13").Value = "numpy.float_" + $nbsp
$ws.Range("C14").Value = "str"
$ws.Range("C15").Value = "numpy.float_" + $nbsp
$ws.Range("C16").Value = "numpy.float_"
$ws.Range("C17").Value = "numpy.float_" + $nbsp
$ws.Range("C18").Value = "numpy.float_" + $nbsp
$ws.Range("C19").Value = "numpy.float_" + $nbsp
$ws.Range("C20").Value = "numpy.float_" + $nbsp
$ws.Range("C21").Value = "numpy.float_" + $nbsp
$ws.Range("C22").Value = "numpy.float_" + $nbsp
$ws.Range("C23").Value = "numpy.float_" + $nbsp
$ws.Range("C24").Value = "numpy.float_" + $nbsp
$ws.Range("C25").Value = "str" + $nbsp
$ws.Range("C26").Value = "str" + $nbsp
$ws.Range("C27").Value = "numpy.float_" + $nbsp
$ws.Range("C28").Value = "numpy.float_" + $nbsp
$ws.Range("C29").Value = "numpy.float_" + $nbsp
$ws.Range("C30").Value = "str"
$ws.Range("C31").Value = "str"

# Move the selected/active cell to B4 (was B14)
$ws.Range("B4").Select() | Out-Null
